$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.200999140739441
$ws.Range("B1").Value = 2.238682746887207
$ws.Range("C1").Value = 3.56810450553894
$ws.Range("D1").Value = 2.577324867248535
$ws.Range("E1").Value = 1.198728680610657
